# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# A new worker record (JOSE ENRIQUE RODRIGUEZ AGUILAR / doc 3876834) is
# inserted as the first data row of the table (row 16). Every worker that
# used to occupy rows 16-37 shifts down one row (17-38) and its "Periodo
# Mora" changes from 2507 to 2508. The row that used to be the extra
# duplicate entry for JOSE ENRIQUE RODRIGUEZ AGUILAR (old row 38, period
# 2505) becomes the new row 38 for the same worker but with period 2508.
#
# Rather than physically moving rows, every cell in the data table
# (B16:G38) is simply (re)written with its final target value - the
# underlying styles for each row are already correct and are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# r | C (doc number) | D (name) | E (periodo) | F (valor mora) | G (salario basico)
$rows = @(
    @(16, "3876834",    "JOSE ENRIQUE RODRIGUEZ AGUILAR",    "2505", 56940, 1423500),
    @(17, "1052958284", "LUIS CARLOS BLANCO WILCHES",        "2508", 56940, 1423500),
    @(18, "33353526",   "DIANA PATRICIA RAMIREZ CARDENAS",   "2508", 56940, 1423500),
    @(19, "1052974162", "ARNOVYS DE JESUS QUEVEDO CARCAMO",  "2508", 56940, 1423500),
    @(20, "55228466",   "KARINA JULIETH ALVARADO CORTISSOZ", "2508", 56940, 1423500),
    @(21, "1052959211", "ALEXANDRA PAOLA PALENCIA VILLALBA", "2508", 56940, 1423500),
    @(22, "1052960424", "LUIS MIGUEL PORTELA MARTINEZ",      "2508", 56940, 1423500),
    @(23, "1052942766", "YULIETH GOMEZ ALDANA",              "2508", 56940, 1423500),
    @(24, "1052997160", "JOSE DAVID GARCIA VILLANUEVA",      "2508", 56940, 1423500),
    @(25, "1053005839", "VICTOR MANUEL ESTRADA GUTIERREZ",   "2508", 56940, 1423500),
    @(26, "1052986392", "CANDELARIA HERNANDEZ ACEVEDO",      "2508", 56940, 1423500),
    @(27, "1140861891", "ALEXANDRA GINNES CAAMAÑO MONTES",   "2508", 56940, 1423500),
    @(28, "1052972068", "MELIZA ANDREA PAVA ACUÑA",          "2508", 56940, 1423500),
    @(29, "1052999956", "JESUS DANIEL FABREGAS CAMACHO",     "2508", 56940, 1423500),
    @(30, "33355060",   "KARINA ELVIRA GARCIA HERNANDEZ",    "2508", 52000, 1300000),
    @(31, "1053003555", "ENELIDA MARIA ORTEGA GARCIA",       "2508", 56940, 1423500),
    @(32, "1002371266", "NAYELIS PAOLA MIELES MERCHAN",      "2508", 56940, 1423500),
    @(33, "1002497507", "WILFRAN PABUENA ROENES",            "2508", 56940, 1423500),
    @(34, "1140825585", "IVANGELA VEGA MUÑOZ",               "2508", 56940, 1423500),
    @(35, "33352604",   "YORJANIS EDITH HERRERA GARCIA",     "2508", 56940, 1423500),
    @(36, "1051743958", "YESENIA JIMENEZ RORIGUEZ",          "2508", 56940, 1423500),
    @(37, "1052951532", "MIRIAM DEL CARMEN MEZA VERGARA",    "2508", 56940, 1423500),
    @(38, "3876834",    "JOSE ENRIQUE RODRIGUEZ AGUILAR",    "2508", 56940, 1423500)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("B$rowNum").Value = "CC"
    $ws.Range("C$rowNum").Value = $r[1]
    $ws.Range("D$rowNum").Value = $r[2]
    $ws.Range("E$rowNum").Value = $r[3]
    $ws.Range("F$rowNum").Value = $r[4]
    $ws.Range("G$rowNum").Value = $r[5]
}
